$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")
$ws.Range("A68").Value = 44049
$ws.Range("B68").Value = 462690
$ws.Range("C68").Value = 506252
$ws.Range("D68").Value = 87973
$ws.Range("E68").Value = 50517
$ws.Range("F68").Value = 26.799152780479368
$ws.Range("G68").Value = 123997
$ws.Range("H68").Value = 9935
$ws.Range("I68").Value = 11928
$ws.Range("J68").Value = 1056915
